$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.704.26'
$ws.Range("E2").Value = '  -0.04%  '

$ws.Range("D3").Value = '2.673.16'
$ws.Range("E3").Value = '  -0.88%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.39'
$ws.Range("E5").Value = '  -1.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.41'
$ws.Range("E6").Value = '  -0.72%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.614'
$ws.Range("E8").Value = '  +4.24%  '

$ws.Range("E9").Value = '  +3.64%  '

$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.399'
$ws.Range("E10").Value = '  -0.51%  '

$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.87'
$ws.Range("E11").Value = '  -1.95%  '

$ws.Range("E12").Value = '  -0.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.26'
$ws.Range("E13").Value = '  -3.63%  '

$ws.Range("E14").Value = '  -3.53%  '

$ws.Range("D15").Value = '3.152.90'
$ws.Range("E15").Value = '  -1.00%  '

$ws.Range("D16").Value = '65.564.71'
$ws.Range("E16").Value = '  -0.09%  '

$ws.Range("D17").Value = '2.668.53'
$ws.Range("E17").Value = '  -0.71%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.89'
$ws.Range("E18").Value = '  +1.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.80'
$ws.Range("E19").Value = '  -2.07%  '

$ws.Range("E20").Value = '  -0.90%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '352.41'
$ws.Range("E21").Value = '  -1.74%  '

$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.73'
$ws.Range("E23").Value = '  -1.84%  '

$ws.Range("E24").Value = '  +5.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.64'
$ws.Range("E25").Value = '  -2.34%  '

$ws.Range("E26").Value = '  +1.55%  '

$ws.Range("E27").Value = '  -2.61%  '

$ws.Range("E28").Value = '  -5.14%  '

$ws.Range("E29").Value = '  -5.76%  '

$ws.Range("E30").Value = '  -0.16%  '

$ws.Range("E31").Value = '  -2.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '529.11'
$ws.Range("E32").Value = '  -3.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.76'
$ws.Range("E33").Value = '  -2.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.54'
$ws.Range("E34").Value = '  +1.51%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.50'
$ws.Range("E35").Value = '  -3.07%  '

$ws.Range("E36").Value = '  -2.52%  '

$ws.Range("E37").Value = '  -1.54%  '

$ws.Range("E38").Value = '  -0.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '157.98'
$ws.Range("E39").Value = '  -3.63%  '

$ws.Range("E40").Value = '  -1.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '163.50'
$ws.Range("E42").Value = '  -4.71%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.12'
$ws.Range("E43").Value = '  -1.74%  '

$ws.Range("E44").Value = '  +3.28%  '

$ws.Range("E45").Value = '  -1.39%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.75'
$ws.Range("E46").Value = '  -3.56%  '

$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₆0264'
$ws.Range("E47").Value = '  +14.89%  '

$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.638'
$ws.Range("E48").Value = '  -2.47%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0257'
$ws.Range("E49").Value = '  -3.37%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.13'
$ws.Range("E50").Value = '  -4.35%  '

$ws.Range("E51").Value = '  +0.13%  '
